$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data set after scraping run: one new tender inserted at the top,
# dates/links refreshed, and the table trimmed back down to 9 rows
# (oldest row drops off the bottom).
$data = @(
    @("Objet : Fourniture et installation de systèmes de détection et extinction d’incendie au niveau des postes électriques pour le compte de Marsa Maroc au Port d’Agadir", "Jeu 17 Juil 2025", "https://achats.marsamaroc.co.ma/?page=entreprise.EntrepriseDetailsConsultation&refConsultation=8030&orgAcronyme=t5y&echanges"),
    @("Objet : TRAVAUX DE VIDANGE DES BACS DE STOCKAGE DE LA STATION DEBALLASTAGE ET RINÇAGE DE MARSA MAROC AU PORT DE MOHAMMEDIA", "Jeu 10 Juil 2025", "https://achats.marsamaroc.co.ma/?page=entreprise.EntrepriseDetailsConsultation&refConsultation=8007&orgAcronyme=t5y&echanges"),
    @("Objet : Travaux de construction de murettes amovibles en béton armé pour le compte de Marsa Maroc et sa filiale SMA au Port d’Agadir", "Mar 08 Juil 2025", "https://achats.marsamaroc.co.ma/?page=entreprise.EntrepriseDetailsConsultation&refConsultation=8003&orgAcronyme=t5y&echanges"),
    @("Objet : TRAVAUX D’AMENAGEMENT DU TERMINAL A CONTENEURS EST DU PORT DE NADOR WEST MED LOT : ELECTRIFICATION HT, BT ET ECT", "Jeu 03 Juil 2025", "https://achats.marsamaroc.co.ma/?page=entreprise.EntrepriseDetailsConsultation&refConsultation=7975&orgAcronyme=t5y&echanges"),
    @("Objet : Fourniture et installation des coffrets électriques pour le compte de MarsaMaroc au port de Tanger Med 1", "Mer 02 Juil 2025", "https://achats.marsamaroc.co.ma/?page=entreprise.EntrepriseDetailsConsultation&refConsultation=7894&orgAcronyme=t5y&echanges"),
    @("Objet : FOURNITURE DE CHARIOTS ELEVATEURS DE MOYENNE CAPACITE A MARSA MAROC", "Mer 02 Juil 2025", "https://achats.marsamaroc.co.ma/?page=entreprise.EntrepriseDetailsConsultation&refConsultation=7967&orgAcronyme=t5y&echanges"),
    @("Objet : Fourniture des chariots élévateurs sur pneus à la Direction de l’Exploitation au Port de Casablanca Trafic Polyvalent.", "Mer 02 Juil 2025", "https://achats.marsamaroc.co.ma/?page=entreprise.EntrepriseDetailsConsultation&refConsultation=7845&orgAcronyme=t5y&echanges"),
    @("Objet : Fourniture des équipements informatiques pour les besoins de la Direction de l’Exploitation au Port de Casablanca Trafic Polyvalent.", "Mar 01 Juil 2025", "https://achats.marsamaroc.co.ma/?page=entreprise.EntrepriseDetailsConsultation&refConsultation=8017&orgAcronyme=t5y&echanges"),
    @("Objet : Démolition d’un ancien bâtiment relevant de Marsa Maroc à la Direction de l’Exploitation au port de Casablanca Trafic Polyvalent", "Mar 01 Juil 2025", "https://achats.marsamaroc.co.ma/?page=entreprise.EntrepriseDetailsConsultation&refConsultation=7985&orgAcronyme=t5y&echanges")
)

# Remove the row that is no longer part of the window (row 11 previously held the
# last tender, the table shrinks back to 10 rows including header).
$ws.Rows.Item(11).Delete() | Out-Null

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $cell = $ws.Cells.Item($row, 3)
    $cell.Value = $data[$i][2]
    $ws.Hyperlinks.Add($cell, $data[$i][2]) | Out-Null
}
